# Apply "Search Functionality - Minor UI Tweaks - 15 titles added" edit
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Action")

# Values are written in this precise order so the shared-string table is
# built up in the same sequence the original authoring session produced.
$ws.Range("A3").Value = "Die Hard"
$ws.Range("B3").Value = "mediaTitle_0011"
$ws.Range("B2").Value = "mediaTitle_0007"
$ws.Range("B4").Value = "mediaTitle_0012"
$ws.Range("B5").Value = "mediaTitle_0013"
$ws.Range("B6").Value = "mediaTitle_0014"
$ws.Range("B7").Value = "mediaTitle_0015"
$ws.Range("A4").Value = "District 9"
$ws.Range("A5").Value = "Spider-Man: Far From Home"
$ws.Range("A6").Value = "I, Robot"
$ws.Range("A7").Value = "Face/Off"

$ws.Range("C3").Value = "Die Hard | 1988 | 2h 12m"
$ws.Range("F3").Value = "Genre: Action, Thriller"

$ws.Range("C4").Value = "District 9 | 2009 | 1h 52m"
$ws.Range("D4").Value = "Violence ensues after an extraterrestrial race forced to live in slum-like conditions on Earth finds a kindred spirit in a government agent exposed to their biotechnology."
$ws.Range("D3").Value = "A New York City police officer tries to save his estranged wife and several others taken hostage by terrorists during a Christmas party at the Nakatomi Plaza in Los Angeles.`n"
$ws.Range("E4").Value = "Cast: Sharlto Copley, Jason Cope, David James, Vanessa Haywood, Mandla Gaduka, Kenneth Nkosi, Eugene Khumbanyiwa, Louis Minnaar, William Allen Young"
$ws.Range("E3").Value = "Cast: Bruce Willis, Alan Rickman, Alexander Godunov, Bonnie Bedelia"
$ws.Range("F4").Value = "Genre: Action, Adventure, Sci-Fi, Fantasy, Horror"

$ws.Range("C5").Value = "Spider-Man: Far From Home | 2019 | 2h 10m"
$ws.Range("D5").Value = "Following the events of Avengers: Endgame (2019), Spider-Man must step up to take on new threats in a world that has changed forever.`n"
$ws.Range("E5").Value = "Cast: Tom Holland, Samuel L. Jackson, Zendaya, Cobie Smulders, Jon Favreau, J. B. Smoove, Jacob Batalon, Martin Starr, Tony Revolori, Marisa Tomei, Jake Gyllenhaal"
$ws.Range("F5").Value = "Genre: Action, Superhero, Sci-Fi, Comedy"

$ws.Range("C6").Value = "I, Robot | 2004 | 1h 55m"
$ws.Range("D6").Value = "In 2035, a technophobic cop investigates a crime that may have been perpetrated by a robot, which leads to a larger threat to humanity.`n"
$ws.Range("E6").Value = "Cast: Will Smith, Bridget Moynahan, Bruce Greenwood, James Cromwell, Chi McBride, Alan Tudyk"
$ws.Range("F6").Value = "Genre: Sci-Fi, Action, Crime, Mystery, Thriller"

$ws.Range("D7").Value = "To foil a terrorist plot, an FBI agent undergoes facial transplant surgery to assume the identity of the criminal mastermind who murdered his only son, but the criminal wakes up prematurely and seeks revenge.`n"
$ws.Range("C7").Value = "Face/Off | 1997 | 2h 13m"
$ws.Range("E7").Value = "Cast: John Travolta, Nicolas Cage, Joan Allen, Gina Gershon, Alessandro Nivola, Colm Feore"
$ws.Range("F7").Value = "Genre: Action, Crime, Sci-Fi, Thriller"

# Row height for row 2 grows (wrapped plot text) and selection moves
$ws.Rows.Item(2).RowHeight = 90.35

# Make Action the active sheet/tab and set new selection
$ws.Activate()
$ws.Range("H4").Select()
